# Updated cryptos list on Fri Oct 25 19:50:21 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: Price values (column D) are always stored as text in this sheet
# (thousands separators use "." so many values are not valid numbers, e.g.
# "67.002.57"). A leading apostrophe forces Excel to keep the assigned
# value as text instead of auto-converting parseable numbers (e.g. "1.00"
# or "585.62") into numeric cells, which would otherwise silently drop
# formatting such as trailing zeros.
$ws.Cells.Item(2, 4).Value = "'67.002.57"
$ws.Cells.Item(2, 5).Value = "  -1.80%  "
$ws.Cells.Item(3, 4).Value = "'2.481.20"
$ws.Cells.Item(3, 5).Value = "  -2.25%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "'585.62"
$ws.Cells.Item(5, 5).Value = "  -1.46%  "
$ws.Cells.Item(6, 4).Value = "'168.22"
$ws.Cells.Item(6, 5).Value = "  -5.04%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 5).Value = "  -3.30%  "
$ws.Cells.Item(9, 4).Value = "'2.481.12"
$ws.Cells.Item(9, 5).Value = "  -2.23%  "
$ws.Cells.Item(10, 5).Value = "  -5.36%  "
$ws.Cells.Item(11, 5).Value = "  +0.15%  "
$ws.Cells.Item(12, 5).Value = "  -4.33%  "
$ws.Cells.Item(13, 5).Value = "  -3.54%  "
$ws.Cells.Item(14, 4).Value = "'25.91"
$ws.Cells.Item(14, 5).Value = "  -4.14%  "
$ws.Cells.Item(15, 4).Value = "'2.936.63"
$ws.Cells.Item(15, 5).Value = "  -2.07%  "
$ws.Cells.Item(16, 5).Value = "  -3.72%  "
$ws.Cells.Item(17, 4).Value = "'66.809.75"
$ws.Cells.Item(17, 5).Value = "  -1.84%  "
$ws.Cells.Item(18, 4).Value = "'2.522.81"
$ws.Cells.Item(18, 5).Value = "  +0.20%  "
$ws.Cells.Item(19, 4).Value = "'11.60"
$ws.Cells.Item(19, 5).Value = "  +0.80%  "
$ws.Cells.Item(20, 4).Value = "'7.76"
$ws.Cells.Item(20, 5).Value = "  -3.46%  "
$ws.Cells.Item(21, 4).Value = "'362.38"
$ws.Cells.Item(21, 5).Value = "  -0.92%  "
$ws.Cells.Item(22, 5).Value = "  -3.93%  "
$ws.Cells.Item(23, 5).Value = "  -5.41%  "
$ws.Cells.Item(24, 5).Value = "  +0.07%  "
$ws.Cells.Item(25, 4).Value = "'70.82"
$ws.Cells.Item(25, 5).Value = "  -0.40%  "
$ws.Cells.Item(26, 5).Value = "  -7.03%  "
$ws.Cells.Item(27, 4).Value = "'9.43"
$ws.Cells.Item(27, 5).Value = "  -8.23%  "
$ws.Cells.Item(28, 4).Value = "'0.999"
$ws.Cells.Item(28, 5).Value = "  +0.15%  "
$ws.Cells.Item(29, 4).Value = "'2.620.50"
$ws.Cells.Item(29, 5).Value = "  -1.72%  "
$ws.Cells.Item(30, 5).Value = "  -7.25%  "
$ws.Cells.Item(31, 5).Value = "  -2.29%  "
$ws.Cells.Item(32, 4).Value = "'515.32"
$ws.Cells.Item(32, 5).Value = "  -6.40%  "
$ws.Cells.Item(33, 5).Value = "  -2.70%  "
$ws.Cells.Item(34, 5).Value = "  -6.51%  "
$ws.Cells.Item(36, 5).Value = "  -3.31%  "
$ws.Cells.Item(37, 4).Value = "'156.78"
$ws.Cells.Item(37, 5).Value = "  +0.08%  "
$ws.Cells.Item(38, 4).Value = "'1.41"
$ws.Cells.Item(39, 4).Value = "'18.92"
$ws.Cells.Item(39, 5).Value = "  +0.05%  "
$ws.Cells.Item(40, 4).Value = "'18.56"
$ws.Cells.Item(40, 5).Value = "  -0.68%  "
$ws.Cells.Item(41, 5).Value = "  -3.88%  "
$ws.Cells.Item(42, 5).Value = "  -5.17%  "
$ws.Cells.Item(43, 4).Value = "'0.333"
$ws.Cells.Item(43, 5).Value = "  -7.07%  "
$ws.Cells.Item(44, 5).Value = "  -3.09%  "
$ws.Cells.Item(45, 4).Value = "'39.18"
$ws.Cells.Item(45, 5).Value = "  -2.33%  "
$ws.Cells.Item(46, 4).Value = "'142.87"
$ws.Cells.Item(46, 5).Value = "  -3.25%  "
$ws.Cells.Item(47, 5).Value = "  -4.92%  "
$ws.Cells.Item(48, 5).Value = "  -3.99%  "
$ws.Cells.Item(49, 5).Value = "  -4.85%  "
$ws.Cells.Item(50, 5).Value = "  -3.58%  "
$ws.Cells.Item(51, 4).Value = "'0.591"
$ws.Cells.Item(51, 5).Value = "  -1.28%  "
